# Updated cryptos list on Mon Aug 21 09:05:12 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'26.206.72"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -0.24%  '

# Row 3
$ws.Cells.Item(3, 5).Value = '  +0.08%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.12%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'216.22"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -0.86%  '

# Row 6
$ws.Cells.Item(6, 4).Value = "'0.5267"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -1.15%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.08%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.32%  '

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.06362"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -1.66%  '

# Row 10
$ws.Cells.Item(10, 4).Value = "'21.43"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -2.28%  '

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.07615"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.98%  '

# Row 12
$ws.Cells.Item(12, 2).Value = 'Polkadot'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(12, 4).Value = "'4.525"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -0.01%  '

# Row 13
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).Value = "'1.672.34"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -0.50%  '

# Row 14
$ws.Cells.Item(14, 4).Value = "'0.5743"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -0.62%  '

# Row 15
$ws.Cells.Item(15, 4).Value = "'0.000008216"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -2.95%  '

# Row 16
$ws.Cells.Item(16, 4).Value = "'66.22"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +2.29%  '

# Row 17
$ws.Cells.Item(17, 4).Value = "'26.212.93"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.46%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  -0.07%  '

# Row 19
$ws.Cells.Item(19, 4).Value = "'4.862"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -0.74%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -1.11%  '

# Row 21
$ws.Cells.Item(21, 4).Value = "'189.85"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.59%  '

# Row 22
$ws.Cells.Item(22, 4).Value = "'6.228"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.35%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -0.09%  '

# Row 24
$ws.Cells.Item(24, 4).Value = "'149.12"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +2.26%  '

# Row 25
$ws.Cells.Item(25, 2).Value = 'Cosmos'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(25, 4).Value = "'7.725"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -1.27%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'Stellar'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(26, 4).Value = "'0.1257"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -1.18%  '

# Row 27
$ws.Cells.Item(27, 4).Value = "'15.89"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +0.75%  '

# Row 28
$ws.Cells.Item(28, 4).Value = "'0.06357"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -2.11%  '

# Row 29
$ws.Cells.Item(29, 4).Value = "'1.375"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -0.43%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -0.43%  '

# Row 31
$ws.Cells.Item(31, 4).Value = "'3.567"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -0.32%  '

# Row 32
$ws.Cells.Item(32, 4).Value = "'3.564"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.53%  '

# Row 33
$ws.Cells.Item(33, 4).Value = "'1.677"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +0.75%  '

# Row 34
$ws.Cells.Item(34, 4).Value = "'1.020"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -1.14%  '

# Row 35
$ws.Cells.Item(35, 4).Value = "'0.6103"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -1.23%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +0.83%  '

# Row 37
$ws.Cells.Item(37, 4).Value = "'2.746"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +1.46%  '

# Row 38
$ws.Cells.Item(38, 4).Value = "'6.177"
$ws.Cells.Item(38, 4).Style = "Normal"

# Row 39
$ws.Cells.Item(39, 5).Value = '  -0.56%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(40, 4).Value = "'0.8871"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +1.84%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'Maker'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(41, 4).Value = "'1.096.62"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -1.36%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -0.43%  '

# Row 43
$ws.Cells.Item(43, 4).Value = "'100.35"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.01%  '

# Row 44
$ws.Cells.Item(44, 4).Value = "'1.830.90"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +0.12%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -0.76%  '

# Row 46
$ws.Cells.Item(46, 4).Value = "'57.42"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +0.52%  '

# Row 47
$ws.Cells.Item(47, 4).Value = "'1.006"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -0.09%  '

# Row 48
$ws.Cells.Item(48, 4).Value = "'8.045"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -1.41%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  +0.11%  '

# Row 50
$ws.Cells.Item(50, 4).Value = "'0.4278"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -0.23%  '

# Row 51
$ws.Cells.Item(51, 4).Value = "'5.995"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -1.40%  '
